$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension to A1:T10 by populating the new 3x3 sending/target-cluster grid
# (ECs x FAPs x sCs) for the Efna2 -> Epha7 LR pair, per Dr Hou's advice.

$rows = @(
  @("ECs", "Efna2", "Epha7", "ECs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 3, 1, 0.3173666666666666, 0.9520999999999999, 0.1863268828340317, 0.1863268828340316, 0.3783271965222221, 3.4049447687, 0.05213314142543136, 0.05213314142543134),
  @("ECs", "Efna2", "Epha7", "FAPs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 2, 0.6666666666666666, 0.2780386666666667, 0.8341160000000001, 0.1632373009158609, 0.1632373009158609, 0.3314449825168889, 2.983004842652, 0.04567281524337266, 0.04567281524337265),
  @("ECs", "Efna2", "Epha7", "sCs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 3, 1, 1.107873666666667, 3.323621, 0.6504358162501074, 0.6504358162501074, 1.320676625598555, 11.886089630387, 0.1819880302883453, 0.1819880302883453),
  @("FAPs", "Efna2", "Epha7", "ECs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 3, 1, 0.3173666666666666, 0.9520999999999999, 0.1863268828340317, 0.1863268828340316, 0.5851484942777777, 5.2663364485, 0.08063292696767507, 0.08063292696767504),
  @("FAPs", "Efna2", "Epha7", "FAPs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 2, 0.6666666666666666, 0.2780386666666667, 0.8341160000000001, 0.1632373009158609, 0.1632373009158609, 0.5126370354511112, 4.613733319060001, 0.07064091430581795, 0.07064091430581794),
  @("FAPs", "Efna2", "Epha7", "sCs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 3, 1, 1.107873666666667, 3.323621, 0.6504358162501074, 0.6504358162501074, 2.042654998109445, 18.383894982985, 0.2814759892461204, 0.2814759892461203),
  @("sCs", "Efna2", "Epha7", "ECs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 3, 1, 0.3173666666666666, 0.9520999999999999, 0.1863268828340317, 0.1863268828340316, 0.3886877371444444, 3.4981896343, 0.05356081444092526, 0.05356081444092525),
  @("sCs", "Efna2", "Epha7", "FAPs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 2, 0.6666666666666666, 0.2780386666666667, 0.8341160000000001, 0.1632373009158609, 0.1632373009158609, 0.3405216474697778, 3.064694827228001, 0.04692357136667033, 0.04692357136667032),
  @("sCs", "Efna2", "Epha7", "sCs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 3, 1, 1.107873666666667, 3.323621, 0.6504358162501074, 0.6504358162501074, 1.356843530738111, 12.211591776643, 0.1869717967156417, 0.1869717967156417)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $rowValues = $rows[$i]
  for ($c = 1; $c -le $rowValues.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
  }
}
